# AFDP-6132 fix user/groups names with correct prefix and domain
#
# The "owning group" literal-participant assignments in the Assignment
# Rules table (Sheet1) referenced LDAP group names that were missing the
# "000." numeric prefix and used the wrong domain (ARMEDIA.COM instead of
# APPDEV.ARMEDIA.COM). Update the affected cells in column G with the
# corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$oldSupervisor = "owning group, ARKCASE_SUPERVISOR@ARMEDIA.COM"
$newSupervisor = "owning group, 000.ARKCASE_SUPERVISOR@APPDEV.ARMEDIA.COM"

$oldEntityAdmin = "owning group, ARKCASE_ENTITY_ADMINISTRATOR@ARMEDIA.COM"
$newEntityAdmin = "owning group, 000.ARKCASE_ENTITY_ADMINISTRATOR@APPDEV.ARMEDIA.COM"

# Organization - Default group / Person - Default group
foreach ($addr in @("G34", "G37")) {
    $cell = $ws.Range($addr)
    if ($cell.Value2 -eq $oldEntityAdmin) {
        $cell.Value = $newEntityAdmin
    }
}

# Complaint - Default group / Case File - Default group / DocumentRepository - Default group
foreach ($addr in @("G23", "G24", "G30")) {
    $cell = $ws.Range($addr)
    if ($cell.Value2 -eq $oldSupervisor) {
        $cell.Value = $newSupervisor
    }
}
